$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.06, 0.04000000000000001, 11.68647933006287, 4.203125),
    @(6.09, 0.04000000000000001, 9.145987272262573, 1.953125),
    @(38.99, 0.16, 59.0127477645874, 18.4375),
    @(39.06, 0.04000000000000001, 12.74627947807312, 4.46875),
    @(39.78, 0.16, 52.52869915962219, 19.015625)
)

$startRow = 249
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
